$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.041452
$ws.Range("H2").Value = 0.124356
$ws.Range("I2").Value = 0.05439747478414846
$ws.Range("J2").Value = 0.05439747478414846
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.572065666666667
$ws.Range("N2").Value = 4.716197
$ws.Range("O2").Value = 0.1759712293834306
$ws.Range("P2").Value = 0.1759712293834305
$ws.Range("Q2").Value = 0.06516526601466667
$ws.Range("R2").Value = 0.586487394132
$ws.Range("S2").Value = 0.00957239051312077
$ws.Range("T2").Value = 0.009572390513120766

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.041452
$ws.Range("H3").Value = 0.124356
$ws.Range("I3").Value = 0.05439747478414846
$ws.Range("J3").Value = 0.05439747478414846
$ws.Range("O3").Value = 0.4743638053196239
$ws.Range("P3").Value = 0.4743638053196239
$ws.Range("Q3").Value = 0.175665327052
$ws.Range("R3").Value = 1.580987943468
$ws.Range("S3").Value = 0.02580419313838695
$ws.Range("T3").Value = 0.02580419313838695

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.041452
$ws.Range("H4").Value = 0.124356
$ws.Range("I4").Value = 0.05439747478414846
$ws.Range("J4").Value = 0.05439747478414846
$ws.Range("M4").Value = 3.123785
$ws.Range("N4").Value = 9.371354999999999
$ws.Range("O4").Value = 0.3496649652969456
$ws.Range("P4").Value = 0.3496649652969455
$ws.Range("Q4").Value = 0.12948713582
$ws.Range("R4").Value = 1.16538422238
$ws.Range("S4").Value = 0.01902089113264074
$ws.Range("T4").Value = 0.01902089113264074

$ws.Range("I5").Value = 0.5204718857143857
$ws.Range("J5").Value = 0.5204718857143856
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.572065666666667
$ws.Range("N5").Value = 4.716197
$ws.Range("O5").Value = 0.1759712293834306
$ws.Range("P5").Value = 0.1759712293834305
$ws.Range("Q5").Value = 0.6234974880785557
$ws.Range("R5").Value = 5.611477392707001
$ws.Range("S5").Value = 0.09158807758867284
$ws.Range("T5").Value = 0.09158807758867279

$ws.Range("I6").Value = 0.5204718857143857
$ws.Range("J6").Value = 0.5204718857143856
$ws.Range("O6").Value = 0.4743638053196239
$ws.Range("P6").Value = 0.4743638053196239
$ws.Range("S6").Value = 0.2468930242693564
$ws.Range("T6").Value = 0.2468930242693564

$ws.Range("I7").Value = 0.5204718857143857
$ws.Range("J7").Value = 0.5204718857143856
$ws.Range("M7").Value = 3.123785
$ws.Range("N7").Value = 9.371354999999999
$ws.Range("O7").Value = 0.3496649652969456
$ws.Range("P7").Value = 0.3496649652969455
$ws.Range("Q7").Value = 1.238925410111667
$ws.Range("R7").Value = 11.150328691005
$ws.Range("S7").Value = 0.1819907838563565
$ws.Range("T7").Value = 0.1819907838563564

$ws.Range("G8").Value = 0.3239583333333333
$ws.Range("H8").Value = 0.971875
$ws.Range("I8").Value = 0.4251306395014658
$ws.Range("J8").Value = 0.4251306395014658
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.572065666666667
$ws.Range("N8").Value = 4.716197
$ws.Range("O8").Value = 0.1759712293834306
$ws.Range("P8").Value = 0.1759712293834305
$ws.Range("Q8").Value = 0.509283773263889
$ws.Range("R8").Value = 4.583553959375
$ws.Range("S8").Value = 0.07481076128163698
$ws.Range("T8").Value = 0.07481076128163695

$ws.Range("G9").Value = 0.3239583333333333
$ws.Range("H9").Value = 0.971875
$ws.Range("I9").Value = 0.4251306395014658
$ws.Range("J9").Value = 0.4251306395014658
$ws.Range("O9").Value = 0.4743638053196239
$ws.Range("P9").Value = 0.4743638053196239
$ws.Range("Q9").Value = 1.372870948958334
$ws.Range("R9").Value = 12.355838540625
$ws.Range("S9").Value = 0.2016665879118806
$ws.Range("T9").Value = 0.2016665879118805

$ws.Range("G10").Value = 0.3239583333333333
$ws.Range("H10").Value = 0.971875
$ws.Range("I10").Value = 0.4251306395014658
$ws.Range("J10").Value = 0.4251306395014658
$ws.Range("M10").Value = 3.123785
$ws.Range("N10").Value = 9.371354999999999
$ws.Range("O10").Value = 0.3496649652969456
$ws.Range("P10").Value = 0.3496649652969455
$ws.Range("Q10").Value = 1.011976182291667
$ws.Range("R10").Value = 9.107785640625
$ws.Range("S10").Value = 0.1486532903079483
$ws.Range("T10").Value = 0.1486532903079483
